# Abrir imagem da galeria - Porem nao salva
# Adds a new "corDark03" / "#010238" row (row 4) to the color table,
# mirroring the existing corDark01/corDark02 rows (A/B text + C color swatch).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the previous row (A3:C3) down onto the new row 4,
# so the new cells reuse the existing cell styles (same style index) instead
# of Excel minting new font/fill/style combinations.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)  # xlPasteFormats

# Label + hex-code cells, mirroring the existing corDark01/corDark02 rows (A2/B2, A3/B3).
# (write B4 before A4 so the new shared-string entries land in the same order as the target file)
$ws.Range("B4").Value = "#010238"
$ws.Range("A4").Value = "corDark03"

# New color swatch cell C4: solid fill #010238, just like C2/C3 are solid swatches
# of corDark01/corDark02.
$ws.Range("C4").ClearContents()
$ws.Range("C4").Interior.Color = 0x01 + (0x02 * 256) + (0x38 * 65536)

# Match the row height used by the other data rows (2 and 3).
$ws.Rows.Item(4).RowHeight = 15.75

# Move the active selection, as captured in the sheetView
$ws.Range("E2").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
